{"js": "// Update the worksheet date and the 25 two-digit multiplication problems.\n// Each old value is unique in the document, so a targeted search +\n// replace keeps every other paragraph/run property untouched.\nconst replacements = [\n  [\"2024-07-15 Monday\", \"2024-07-16 Tuesday\"],\n  [\"76\\u00d786=\", \"72\\u00d770=\"],\n  [\"95\\u00d766=\", \"57\\u00d789=\"],\n  [\"53\\u00d769=\", \"60\\u00d756=\"],\n  [\"67\\u00d790=\", \"45\\u00d737=\"],\n  [\"35\\u00d770=\", \"68\\u00d761=\"],\n  [\"54\\u00d767=\", \"67\\u00d723=\"],\n  [\"53\\u00d785=\", \"71\\u00d721=\"],\n  [\"95\\u00d729=\", \"33\\u00d775=\"],\n  [\"26\\u00d759=\", \"37\\u00d732=\"],\n  [\"79\\u00d790=\", \"76\\u00d775=\"],\n  [\"38\\u00d738=\", \"97\\u00d781=\"],\n  [\"36\\u00d783=\", \"67\\u00d760=\"],\n  [\"20\\u00d747=\", \"86\\u00d715=\"],\n  [\"14\\u00d797=\", \"90\\u00d718=\"],\n  [\"86\\u00d762=\", \"55\\u00d773=\"],\n  [\"82\\u00d723=\", \"71\\u00d791=\"],\n  [\"97\\u00d768=\", \"82\\u00d750=\"],\n  [\"90\\u00d797=\", \"23\\u00d787=\"],\n  [\"83\\u00d781=\", \"48\\u00d736=\"],\n  [\"59\\u00d750=\", \"71\\u00d759=\"],\n  [\"66\\u00d786=\", \"38\\u00d719=\"],\n  [\"97\\u00d748=\", \"53\\u00d726=\"],\n  [\"78\\u00d744=\", \"51\\u00d749=\"],\n  [\"93\\u00d717=\", \"93\\u00d795=\"],\n  [\"86\\u00d712=\", \"20\\u00d751=\"],\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  for (const range of results.items) {\n    range.insertText(newText, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n", "ps1": "# Update the worksheet date and the 25 two-digit multiplication problems.\n# Each old value is unique in the document, so Find/Replace against the\n# whole document body (Content) swaps only the target w:t text while\n# leaving paragraph/run formatting untouched.\n$d = $word.ActiveDocument\n\n$replacements = @(\n    @(\"2024-07-15 Monday\", \"2024-07-16 Tuesday\"),\n    @(\"76\u00d786=\", \"72\u00d770=\"),\n    @(\"95\u00d766=\", \"57\u00d789=\"),\n    @(\"53\u00d769=\", \"60\u00d756=\"),\n    @(\"67\u00d790=\", \"45\u00d737=\"),\n    @(\"35\u00d770=\", \"68\u00d761=\"),\n    @(\"54\u00d767=\", \"67\u00d723=\"),\n    @(\"53\u00d785=\", \"71\u00d721=\"),\n    @(\"95\u00d729=\", \"33\u00d775=\"),\n    @(\"26\u00d759=\", \"37\u00d732=\"),\n    @(\"79\u00d790=\", \"76\u00d775=\"),\n    @(\"38\u00d738=\", \"97\u00d781=\"),\n    @(\"36\u00d783=\", \"67\u00d760=\"),\n    @(\"20\u00d747=\", \"86\u00d715=\"),\n    @(\"14\u00d797=\", \"90\u00d718=\"),\n    @(\"86\u00d762=\", \"55\u00d773=\"),\n    @(\"82\u00d723=\", \"71\u00d791=\"),\n    @(\"97\u00d768=\", \"82\u00d750=\"),\n    @(\"90\u00d797=\", \"23\u00d787=\"),\n    @(\"83\u00d781=\", \"48\u00d736=\"),\n    @(\"59\u00d750=\", \"71\u00d759=\"),\n    @(\"66\u00d786=\", \"38\u00d719=\"),\n    @(\"97\u00d748=\", \"53\u00d726=\"),\n    @(\"78\u00d744=\", \"51\u00d749=\"),\n    @(\"93\u00d717=\", \"93\u00d795=\"),\n    @(\"86\u00d712=\", \"20\u00d751=\")\n)\n\nforeach ($pair in $replacements) {\n    $oldText = $pair[0]\n    $newText = $pair[1]\n\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Execute($oldText, $true, $false, $false, $false, $false, $true, 1, $false, $newText, 2)\n}\n"}
